$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.851.65'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.68%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.966.20'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '595.63'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.33%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.92'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.50%  '

$ws.Range("E7").Value = '  -0.11%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '2.965.86'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.11%  '

$ws.Range("E9").Value = '  +0.14%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.31'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.36%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.149'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.25%  '

$ws.Range("E12").Value = '  +0.75%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000240'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +6.25%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.12'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.08%  '

$ws.Range("E15").Value = '  -0.54%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.456.22'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.02%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.695.90'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.59%  '

$ws.Range("E18").Value = '  -0.49%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.967.06'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.90%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '442.27'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.88%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.50'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.03%  '

$ws.Range("E22").Value = '  -1.35%  '

$ws.Range("E23").Value = '  -0.69%  '

$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '81.59'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.34%  '

$ws.Range("B25").Value = 'RenderToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.25'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.78%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.88'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.32%  '

$ws.Range("E27").Value = '  -3.67%  '

$ws.Range("E28").Value = '  -0.03%  '

$ws.Range("E29").Value = '  +3.90%  '

$ws.Range("E30").Value = '  -0.32%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.15'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.36%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0₃0960'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +9.83%  '

$ws.Range("E33").Value = '  -0.66%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '26.53'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.90%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.06%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.993'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.01%  '

$ws.Range("E37").Value = '  -0.14%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.06'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.43%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.04'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.01%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '49.48'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.00%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.54'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.01%  '

$ws.Range("E42").Value = '  -5.27%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.281'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.47%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '40.43'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.01%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.741.14'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.27%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '134.97'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.06%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0339'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.31%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '363.95'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.40%  '

$ws.Range("E49").Value = '  +0.05%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '23.04'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.93%  '

$ws.Range("E51").Value = '  -0.77%  '
